# STEP 8 + STEP 9. Final.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# STEP 8 — rewrite the per-flight "xx" duration column (E3:E9)
$ws.Range("E3").Value = 144
$ws.Range("E4").Value = 144
$ws.Range("E5").Value = 216
$ws.Range("E6").Value = 72
$ws.Range("E7").Value = 144
$ws.Range("E8").Value = 72
$ws.Range("E9").Value = 144

# STEP 9 — rewrite the "sum" column (E11:E17)
$ws.Range("E11").Value = 17280
$ws.Range("E12").Value = 17280
$ws.Range("E13").Value = 25920
$ws.Range("E14").Value = 8640
$ws.Range("E15").Value = 17280
$ws.Range("E16").Value = 8640
$ws.Range("E17").Value = 17280

# Leave the cursor parked on the (now unused) column F, matching the
# author's final selection state in the saved workbook.
$ws.Columns("F").Select()
